$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -3
$ws.Range("F3").Value = -7
$ws.Range("F4").Value = -5
$ws.Range("F5").Value = 4
$ws.Range("F7").Value = -1
$ws.Range("F9").Value = -2
